# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Integral"     (used by the slide master / main presentation)
#   ppt/theme/theme2.xml  -> "Office Theme" (used by the notes master)
#
# The commit swaps their colour schemes: theme1.xml ends up carrying the
# stock "Office" palette while theme2.xml ends up carrying the old
# "Integral" palette (font scheme / format scheme are identical between the
# two parts already, so only the 12 clrScheme colours actually move).
#
# The PowerPoint object model only exposes a writable colour scheme on the
# slide/master side (Slide.ThemeColorScheme / Master.ColorScheme), which is
# backed by theme1.xml, so that is the half of the swap we can perform here.
# ThemeColorScheme is the modern 12-slot scheme (dk1,lt1,dk2,lt2,accent1-6,
# hlink,folHlink) - unlike the legacy 8-slot ColorScheme object it also
# leaves accent5/accent6/hlink/folHlink reachable and does not blow away the
# <a:clrScheme> name attribute.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() isn't available in this host, so the values below are the
# 0x00BBGGRR COM colour encoding of the target "Office" theme hex colours:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$tcs.Colors(1).RGB = 0          # dk1      -> 000000
$tcs.Colors(2).RGB = 16777215   # lt1      -> FFFFFF
$tcs.Colors(3).RGB = 6968388    # dk2      -> 44546A
$tcs.Colors(4).RGB = 15132391   # lt2      -> E7E6E6
$tcs.Colors(5).RGB = 13998939   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB = 3243501    # accent2  -> ED7D31
$tcs.Colors(7).RGB = 10855845   # accent3  -> A5A5A5
$tcs.Colors(8).RGB = 49407      # accent4  -> FFC000
$tcs.Colors(9).RGB = 12874308   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    -> 0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink -> 954F72
